$wb = $excel.ActiveWorkbook

# --- Sheet1: remove the now-unused trailing rows (45:87) that only had
#     leftover index numbers in column A, and update the selection/view ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows("45:87").Delete()

# --- Sheet3: it is no longer the tab that is active/selected ---
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Make Sheet1 the active sheet/tab again, with the view scrolled down
#     and F61 as the active cell (matches the author's final view state) ---
$ws1.Activate()
$ws1.Range("F61").Select()
